$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 867-868, shifting existing rows 867-925 down to 869-927
$ws.Range("A867:A868").EntireRow.Insert()

# Populate new row 867 (Primera)
$ws.Cells.Item(867,1).Value = 6
$ws.Cells.Item(867,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(867,3).Value = "Metropolitana"
$ws.Cells.Item(867,4).Value = 44826
$ws.Cells.Item(867,5).Value = 13
$ws.Cells.Item(867,6).Value = 100112023
$ws.Cells.Item(867,7).Value = "Brócoli"
$ws.Cells.Item(867,8).Value = "Sin especificar"
$ws.Cells.Item(867,9).Value = "Primera"
$ws.Cells.Item(867,10).Value = 8400
$ws.Cells.Item(867,11).Value = 750
$ws.Cells.Item(867,12).Value = 800
$ws.Cells.Item(867,13).Value = 771
$ws.Cells.Item(867,14).Value = '$/unidad'
$ws.Cells.Item(867,15).Value = "Región Metropolitana"
$ws.Cells.Item(867,16).Value = 771
$ws.Cells.Item(867,17).Value = 1
$ws.Cells.Item(867,18).Value = "Hortaliza"

# Populate new row 868 (Segunda)
$ws.Cells.Item(868,1).Value = 6
$ws.Cells.Item(868,2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(868,3).Value = "Metropolitana"
$ws.Cells.Item(868,4).Value = 44826
$ws.Cells.Item(868,5).Value = 13
$ws.Cells.Item(868,6).Value = 100112023
$ws.Cells.Item(868,7).Value = "Brócoli"
$ws.Cells.Item(868,8).Value = "Sin especificar"
$ws.Cells.Item(868,9).Value = "Segunda"
$ws.Cells.Item(868,10).Value = 3100
$ws.Cells.Item(868,11).Value = 650
$ws.Cells.Item(868,12).Value = 650
$ws.Cells.Item(868,13).Value = 650
$ws.Cells.Item(868,14).Value = '$/unidad'
$ws.Cells.Item(868,15).Value = "Región Metropolitana"
$ws.Cells.Item(868,16).Value = 650
$ws.Cells.Item(868,17).Value = 1
$ws.Cells.Item(868,18).Value = "Hortaliza"

Write-Output "done"
